# Update the "Förändrad" (Changed) date in column C for all data rows
# (rows 2-181) from 2023-09-11 (45180) to 2023-09-12 (45181).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C181").Value = 45181
